{"js": "// Applies the three text-level edits described by the diff:\n//   1. \"Running the Reactive Agent\" paragraph: drop \"compiled and\", replace\n//      \"the demo agent (i.e., \"javac Krislet.java\" and \"java X\")\" with\n//      \"Krislet (i.e., \"java Krislet\")\".\n//   2. \"Where <ball visibility> ...\" paragraph: add \"DirectlyInFront},\" to\n//      the <goal visibility> set (matching the <ball visibility> set above it).\n//   3. \"Expected Action of the Provided Agent\" paragraph: \"the demo agent\"\n//      becomes \"Krislet\".\n\nasync function replaceOnce(context, searchText, replacement) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Running the Reactive Agent\" paragraph.\nawait replaceOnce(\n  context,\n  \"The agent is compiled and executed identically to the demo agent (i.e., \\u201cjavac Krislet.java\\u201d and \\u201cjava Krislet\\u201d).\",\n  \"The agent is executed identically to Krislet (i.e., \\u201cjava Krislet\\u201d).\"\n);\n\n// 2) Add \"DirectlyInFront\" to the <goal visibility> set.\nawait replaceOnce(\n  context,\n  \"\\u2208 {Visible, NotVisible, and <action>\",\n  \"\\u2208 {Visible, NotVisible, DirectlyInFront},and <action>\"\n);\n\n// 3) \"Expected Action of the Provided Agent\" paragraph.\nawait replaceOnce(\n  context,\n  \"is intended to behave similarly to the demo agent. The agent will turn\",\n  \"is intended to behave similarly to Krislet. The agent will turn\"\n);\n", "ps1": "# Applies the three text-level edits described by the diff:\n#   1. \"Running the Reactive Agent\" paragraph: drop \"compiled and\", replace\n#      \"the demo agent (i.e., \"javac Krislet.java\" and \"java X\")\" with\n#      \"Krislet (i.e., \"java Krislet\")\".\n#   2. \"Where <ball visibility> ...\" paragraph: add \"DirectlyInFront},\" to\n#      the <goal visibility> set (matching the <ball visibility> set above it).\n#   3. \"Expected Action of the Provided Agent\" paragraph: \"the demo agent\"\n#      becomes \"Krislet\".\n#\n# NOTE: the Word enum constants (wdReplaceOne, wdFindContinue, ...) are NOT\n# predefined globals in this host, so the literal values are used directly:\n#   wdFindContinue = 1, wdReplaceOne = 1\n# NOTE: the \"\u2208\" (U+2208) math glyph living inside the <m:oMath> runs is not\n# matchable by Find.Execute in this host, so the searches below are anchored\n# on the plain-text runs immediately around it instead of including it.\n\n$d = $word.ActiveDocument\n\n# 1) \"Running the Reactive Agent\" paragraph.\n$r1 = $d.Content\n$r1.Find.Execute(\n    \"The agent is compiled and executed identically to the demo agent (i.e., \u201cjavac Krislet.java\u201d and \u201cjava Krislet\u201d).\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"The agent is executed identically to Krislet (i.e., \u201cjava Krislet\u201d).\",\n    1\n)\n\n# 2) Add \"DirectlyInFront\" to the <goal visibility> set.\n$r2 = $d.Content\n$r2.Find.Execute(\n    \"{Visible, NotVisible, and <action>\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"{Visible, NotVisible, DirectlyInFront},and <action>\",\n    1\n)\n\n# 3) \"Expected Action of the Provided Agent\" paragraph.\n$r3 = $d.Content\n$r3.Find.Execute(\n    \"is intended to behave similarly to the demo agent. The agent will turn\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"is intended to behave similarly to Krislet. The agent will turn\",\n    1\n)\n"}
